$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 7 (shifts existing rows 7+ down by one)
$null = $ws.Rows("7:7").Insert()

# Populate the new row 7
$ws.Range("B7").Value = "ngx-timeago"
$ws.Range("C7").Value = "Live updating timestamps in Angular 6+."

# Update the active cell selection
$null = $ws.Range("C8").Select()
